# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to the Leve profit tables
# across all job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 456.125
$ws.Range("I28").Value = 451.2857
$ws.Range("J28").Value = 459.8889
$ws.Range("K28").Value = 451.2857
$ws.Range("L28").Value = 459.8889
$ws.Range("M28").Value = 33.71429999999998
$ws.Range("N28").Value = -1429.8889
$ws.Range("H40").Value = 6406.7144
$ws.Range("I40").Value = 6613.05
$ws.Range("K40").Value = 6613.05
$ws.Range("M40").Value = -6438.05
$ws.Range("H112").Value = 1452.4
$ws.Range("I112").Value = 379.83334
$ws.Range("J112").Value = 1641.6765
$ws.Range("K112").Value = 1139.50002
$ws.Range("L112").Value = 4925.029500000001
$ws.Range("M112").Value = -31.50001999999995
$ws.Range("N112").Value = -7141.029500000001
$ws.Range("H132").Value = 26706.13
$ws.Range("I132").Value = 4057.5
$ws.Range("J132").Value = 75778.164
$ws.Range("K132").Value = 12172.5
$ws.Range("L132").Value = 227334.492
$ws.Range("M132").Value = -9642.5
$ws.Range("N132").Value = -232394.492
$ws.Range("H138").Value = 2724.9636
$ws.Range("I138").Value = 1787.9584
$ws.Range("J138").Value = 3450.3872
$ws.Range("K138").Value = 5363.8752
$ws.Range("L138").Value = 10351.1616
$ws.Range("M138").Value = -223.8752000000004
$ws.Range("N138").Value = -20631.1616

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10778.42
$ws.Range("I32").Value = 9919.164000000001
$ws.Range("J32").Value = 17330.25
$ws.Range("K32").Value = 9919.164000000001
$ws.Range("L32").Value = 17330.25
$ws.Range("M32").Value = -9632.164000000001
$ws.Range("N32").Value = -17904.25
$ws.Range("H61").Value = 1951.1765
$ws.Range("I61").Value = 1141.9474
$ws.Range("J61").Value = 2976.2
$ws.Range("K61").Value = 1141.9474
$ws.Range("L61").Value = 2976.2
$ws.Range("M61").Value = -929.9474
$ws.Range("N61").Value = -3400.2
$ws.Range("H74").Value = 1936.7858
$ws.Range("I74").Value = 1529.4857
$ws.Range("J74").Value = 3973.2856
$ws.Range("K74").Value = 1529.4857
$ws.Range("L74").Value = 3973.2856
$ws.Range("M74").Value = -655.4857
$ws.Range("N74").Value = -5721.2856
$ws.Range("H77").Value = 1936.7858
$ws.Range("I77").Value = 1529.4857
$ws.Range("J77").Value = 3973.2856
$ws.Range("K77").Value = 7647.4285
$ws.Range("L77").Value = 19866.428
$ws.Range("M77").Value = -3279.4285
$ws.Range("N77").Value = -28602.428
$ws.Range("H122").Value = 2084.6667
$ws.Range("I122").Value = 1725.25
$ws.Range("J122").Value = 2803.5
$ws.Range("K122").Value = 5175.75
$ws.Range("L122").Value = 8410.5
$ws.Range("M122").Value = -2725.75
$ws.Range("N122").Value = -13310.5
$ws.Range("H136").Value = 1951.1765
$ws.Range("I136").Value = 1141.9474
$ws.Range("J136").Value = 2976.2
$ws.Range("K136").Value = 3425.8422
$ws.Range("L136").Value = 8928.599999999999
$ws.Range("M136").Value = -875.8422
$ws.Range("N136").Value = -14028.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1167
$ws.Range("I22").Value = 1167
$ws.Range("K22").Value = 1167
$ws.Range("M22").Value = -994
$ws.Range("H86").Value = 2251.6667
$ws.Range("I86").Value = 1852.2858
$ws.Range("J86").Value = 3649.5
$ws.Range("K86").Value = 1852.2858
$ws.Range("L86").Value = 3649.5
$ws.Range("M86").Value = -729.2858000000001
$ws.Range("N86").Value = -5895.5
$ws.Range("H89").Value = 2251.6667
$ws.Range("I89").Value = 1852.2858
$ws.Range("J89").Value = 3649.5
$ws.Range("K89").Value = 9261.429
$ws.Range("L89").Value = 18247.5
$ws.Range("M89").Value = -3645.429
$ws.Range("N89").Value = -29479.5
$ws.Range("H99").Value = 2436.3333
$ws.Range("I99").Value = 2190
$ws.Range("J99").Value = 2929
$ws.Range("K99").Value = 2190
$ws.Range("L99").Value = 2929
$ws.Range("M99").Value = -692
$ws.Range("N99").Value = -5925

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H137").Value = 63795.6
$ws.Range("J137").Value = 63795.6
$ws.Range("L137").Value = 63795.6
$ws.Range("N137").Value = -73995.60000000001
$ws.Range("H138").Value = 40689.168
$ws.Range("J138").Value = 40689.168
$ws.Range("L138").Value = 40689.168
$ws.Range("N138").Value = -50969.168
$ws.Range("H139").Value = 100890
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 100890
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 100890
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -111170
$ws.Range("H141").Value = 12977
$ws.Range("J141").Value = 12977
$ws.Range("L141").Value = 12977
$ws.Range("N141").Value = -23337

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 384.77777
$ws.Range("I26").Value = 374.5
$ws.Range("J26").Value = 405.33334
$ws.Range("K26").Value = 1123.5
$ws.Range("L26").Value = 1216.00002
$ws.Range("M26").Value = -835.5
$ws.Range("N26").Value = -1792.00002
$ws.Range("H68").Value = 1340.4482
$ws.Range("I68").Value = 995.4761999999999
$ws.Range("J68").Value = 1450.2122
$ws.Range("K68").Value = 2986.4286
$ws.Range("L68").Value = 4350.6366
$ws.Range("M68").Value = -2175.4286
$ws.Range("N68").Value = -5972.6366
$ws.Range("H71").Value = 1340.4482
$ws.Range("I71").Value = 995.4761999999999
$ws.Range("J71").Value = 1450.2122
$ws.Range("K71").Value = 8959.2858
$ws.Range("L71").Value = 13051.9098
$ws.Range("M71").Value = -4903.2858
$ws.Range("N71").Value = -21163.9098
$ws.Range("H107").Value = 9791.739
$ws.Range("I107").Value = 9516.272000000001
$ws.Range("J107").Value = 10044.25
$ws.Range("K107").Value = 28548.816
$ws.Range("L107").Value = 30132.75
$ws.Range("M107").Value = -26628.816
$ws.Range("N107").Value = -33972.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I25").Value = 10000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 10000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -9471
$ws.Range("N25").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 1389.7778
$ws.Range("I122").Value = 1380
$ws.Range("J122").Value = 1402
$ws.Range("K122").Value = 4140
$ws.Range("L122").Value = 4206
$ws.Range("M122").Value = -1690
$ws.Range("N122").Value = -9106
$ws.Range("H132").Value = 47624590
$ws.Range("I132").Value = 66672280
$ws.Range("K132").Value = 200016840
$ws.Range("M132").Value = -200014310
$ws.Range("H133").Value = 61996
$ws.Range("J133").Value = 61996
$ws.Range("L133").Value = 61996
$ws.Range("N133").Value = -72116

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 21542.5
$ws.Range("J63").Value = 21542.5
$ws.Range("L63").Value = 21542.5
$ws.Range("N63").Value = -23040.5
$ws.Range("H66").Value = 21542.5
$ws.Range("J66").Value = 21542.5
$ws.Range("L66").Value = 64627.5
$ws.Range("N66").Value = -72115.5
$ws.Range("H122").Value = 39309.742
$ws.Range("I122").Value = 49840.617
$ws.Range("K122").Value = 149521.851
$ws.Range("M122").Value = -147071.851
$ws.Range("H124").Value = 42714.5
$ws.Range("J124").Value = 42714.5
$ws.Range("L124").Value = 42714.5
$ws.Range("N124").Value = -52534.5
$ws.Range("H132").Value = 3872.1177
$ws.Range("I132").Value = 3247.353
$ws.Range("J132").Value = 4496.8823
$ws.Range("K132").Value = 9742.059000000001
$ws.Range("L132").Value = 13490.6469
$ws.Range("M132").Value = -7212.059000000001
$ws.Range("N132").Value = -18550.6469

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1813847.5
$ws.Range("I132").Value = 2558991
$ws.Range("J132").Value = 4213
$ws.Range("K132").Value = 7676973
$ws.Range("L132").Value = 12639
$ws.Range("M132").Value = -7674443
$ws.Range("N132").Value = -17699
$ws.Range("H136").Value = 389758.8
$ws.Range("I136").Value = 556148.1
$ws.Range("J136").Value = 1517
$ws.Range("K136").Value = 1668444.3
$ws.Range("L136").Value = 4551
$ws.Range("M136").Value = -1665894.3
$ws.Range("N136").Value = -9651
